$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Cells.Item(5, 8).Value = 580
$ws.Cells.Item(5, 9).Value = 655.2308
$ws.Cells.Item(5, 11).Value = 655.2308
$ws.Cells.Item(5, 13).Value = -540.2308
# Row 47
$ws.Cells.Item(47, 8).Value = 47000
$ws.Cells.Item(47, 9).Value = 47000
$ws.Cells.Item(47, 11).Value = 47000
$ws.Cells.Item(47, 13).Value = -46028
# Row 74
$ws.Cells.Item(74, 8).Value = 4728.7144
$ws.Cells.Item(74, 9).Value = 3614.7144
$ws.Cells.Item(74, 10).Value = 5285.7144
$ws.Cells.Item(74, 11).Value = 3614.7144
$ws.Cells.Item(74, 12).Value = 5285.7144
$ws.Cells.Item(74, 13).Value = -2678.7144
$ws.Cells.Item(74, 14).Value = -7157.7144
# Row 77
$ws.Cells.Item(77, 8).Value = 4728.7144
$ws.Cells.Item(77, 9).Value = 3614.7144
$ws.Cells.Item(77, 10).Value = 5285.7144
$ws.Cells.Item(77, 11).Value = 18073.572
$ws.Cells.Item(77, 12).Value = 26428.572
$ws.Cells.Item(77, 13).Value = -13393.572
$ws.Cells.Item(77, 14).Value = -35788.572
# Row 86
$ws.Cells.Item(86, 8).Value = 4183.8096
$ws.Cells.Item(86, 10).Value = 4954.778
$ws.Cells.Item(86, 12).Value = 4954.778
$ws.Cells.Item(86, 14).Value = -7200.778
# Row 89
$ws.Cells.Item(89, 8).Value = 4183.8096
$ws.Cells.Item(89, 10).Value = 4954.778
$ws.Cells.Item(89, 12).Value = 24773.89
$ws.Cells.Item(89, 14).Value = -36005.89
# Row 106
$ws.Cells.Item(106, 8).Value = 2297.4375
$ws.Cells.Item(106, 9).Value = 2054.2144
$ws.Cells.Item(106, 11).Value = 2054.2144
$ws.Cells.Item(106, 13).Value = -1423.2144
# Row 107
$ws.Cells.Item(107, 8).Value = 1352.1177
$ws.Cells.Item(107, 9).Value = 821.63635
$ws.Cells.Item(107, 11).Value = 821.63635
$ws.Cells.Item(107, 13).Value = 1098.36365
# Row 111
$ws.Cells.Item(111, 8).Value = 4847.095
$ws.Cells.Item(111, 9).Value = 4926.8667
$ws.Cells.Item(111, 11).Value = 14780.6001
$ws.Cells.Item(111, 13).Value = -11713.6001
# Row 125
$ws.Cells.Item(125, 8).Value = 5739.8
$ws.Cells.Item(125, 9).Value = 4200
$ws.Cells.Item(125, 10).Value = 6124.75
$ws.Cells.Item(125, 11).Value = 37800
$ws.Cells.Item(125, 12).Value = 55122.75
$ws.Cells.Item(125, 13).Value = -35340
$ws.Cells.Item(125, 14).Value = -60042.75
# Row 132
$ws.Cells.Item(132, 8).Value = 3533.9
$ws.Cells.Item(132, 9).Value = 3238.8333
$ws.Cells.Item(132, 11).Value = 9716.499899999999
$ws.Cells.Item(132, 13).Value = -7186.499899999999
# Row 137
$ws.Cells.Item(137, 8).Value = 2411.2307
$ws.Cells.Item(137, 9).Value = 2471.2942
$ws.Cells.Item(137, 11).Value = 7413.882599999999
$ws.Cells.Item(137, 13).Value = -4863.882599999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 2637.75
$ws.Cells.Item(45, 9).Value = 2212
$ws.Cells.Item(45, 11).Value = 2212
$ws.Cells.Item(45, 13).Value = -1835
# Row 74
$ws.Cells.Item(74, 8).Value = 166670850
$ws.Cells.Item(74, 9).Value = 333336640
$ws.Cells.Item(74, 11).Value = 333336640
$ws.Cells.Item(74, 13).Value = -333335766
# Row 77
$ws.Cells.Item(77, 8).Value = 166670850
$ws.Cells.Item(77, 9).Value = 333336640
$ws.Cells.Item(77, 11).Value = 1666683200
$ws.Cells.Item(77, 13).Value = -1666678832
# Row 102
$ws.Cells.Item(102, 8).Value = 2240
$ws.Cells.Item(102, 9).Value = 1614
$ws.Cells.Item(102, 11).Value = 1614
$ws.Cells.Item(102, 13).Value = 8
# Row 110
$ws.Cells.Item(110, 8).Value = 3605.2222
$ws.Cells.Item(110, 9).Value = 2633.3076
$ws.Cells.Item(110, 11).Value = 2633.3076
$ws.Cells.Item(110, 13).Value = -588.3076000000001
# Row 122
$ws.Cells.Item(122, 8).Value = 3528.6206
$ws.Cells.Item(122, 9).Value = 2342.7144
$ws.Cells.Item(122, 11).Value = 7028.1432
$ws.Cells.Item(122, 13).Value = -4578.1432
# Row 132
$ws.Cells.Item(132, 8).Value = 47621896
$ws.Cells.Item(132, 9).Value = 2877.5557
$ws.Cells.Item(132, 10).Value = 333336000
$ws.Cells.Item(132, 11).Value = 8632.667099999999
$ws.Cells.Item(132, 12).Value = 1000008000
$ws.Cells.Item(132, 13).Value = -6102.667099999999
$ws.Cells.Item(132, 14).Value = -1000013060

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 2184.8572
$ws.Cells.Item(20, 9).Value = 2557
$ws.Cells.Item(20, 10).Value = 1812.7142
$ws.Cells.Item(20, 11).Value = 2557
$ws.Cells.Item(20, 12).Value = 1812.7142
$ws.Cells.Item(20, 13).Value = -2310
$ws.Cells.Item(20, 14).Value = -2306.7142
# Row 86
$ws.Cells.Item(86, 8).Value = 10522.571
$ws.Cells.Item(86, 9).Value = 7623.8125
$ws.Cells.Item(86, 11).Value = 7623.8125
$ws.Cells.Item(86, 13).Value = -6500.8125
# Row 89
$ws.Cells.Item(89, 8).Value = 10522.571
$ws.Cells.Item(89, 9).Value = 7623.8125
$ws.Cells.Item(89, 11).Value = 38119.0625
$ws.Cells.Item(89, 13).Value = -32503.0625
# Row 105
$ws.Cells.Item(105, 8).Value = 6132.905
$ws.Cells.Item(105, 9).Value = 9408
$ws.Cells.Item(105, 11).Value = 9408
$ws.Cells.Item(105, 13).Value = -7661

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Cells.Item(5, 8).Value = 648.2222
$ws.Cells.Item(5, 10).Value = 247.33333
$ws.Cells.Item(5, 12).Value = 247.33333
$ws.Cells.Item(5, 14).Value = -471.33333
# Row 31
$ws.Cells.Item(31, 8).Value = 29417400
$ws.Cells.Item(31, 9).Value = 4100.9614
$ws.Cells.Item(31, 10).Value = 125010620
$ws.Cells.Item(31, 11).Value = 4100.9614
$ws.Cells.Item(31, 12).Value = 125010620
$ws.Cells.Item(31, 13).Value = -3805.9614
$ws.Cells.Item(31, 14).Value = -125011210
# Row 34
$ws.Cells.Item(34, 8).Value = 29417400
$ws.Cells.Item(34, 9).Value = 4100.9614
$ws.Cells.Item(34, 10).Value = 125010620
$ws.Cells.Item(34, 11).Value = 4100.9614
$ws.Cells.Item(34, 12).Value = 125010620
$ws.Cells.Item(34, 13).Value = -3898.9614
$ws.Cells.Item(34, 14).Value = -125011024
# Row 131
$ws.Cells.Item(131, 8).Value = 68180.836
$ws.Cells.Item(131, 9).Value = 35000
$ws.Cells.Item(131, 11).Value = 35000
$ws.Cells.Item(131, 13).Value = -29960
# Row 134
$ws.Cells.Item(134, 8).Value = 1291.8462
$ws.Cells.Item(134, 9).Value = 1291.8462
$ws.Cells.Item(134, 11).Value = 3875.5386
$ws.Cells.Item(134, 13).Value = -1340.5386

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Cells.Item(4, 8).Value = 3116270.2
$ws.Cells.Item(4, 9).Value = 1875438.6
$ws.Cells.Item(4, 10).Value = 8700012
$ws.Cells.Item(4, 11).Value = 5626315.800000001
$ws.Cells.Item(4, 12).Value = 26100036
$ws.Cells.Item(4, 13).Value = -5626203.800000001
$ws.Cells.Item(4, 14).Value = -26100260
# Row 49
$ws.Cells.Item(49, 8).Value = 10000
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 10000
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 30000
$ws.Cells.Item(49, 13).ClearContents()
$ws.Cells.Item(49, 14).Value = -30312
# Row 69
$ws.Cells.Item(69, 8).Value = 5000
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 5000
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 15000
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).Value = -16622
# Row 72
$ws.Cells.Item(72, 8).Value = 5000
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 5000
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 45000
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).Value = -53112
# Row 113
$ws.Cells.Item(113, 8).Value = 2765.6316
$ws.Cells.Item(113, 9).Value = 1732.7778
$ws.Cells.Item(113, 10).Value = 3695.2
$ws.Cells.Item(113, 11).Value = 5198.3334
$ws.Cells.Item(113, 12).Value = 11085.6
$ws.Cells.Item(113, 13).Value = -3028.3334
$ws.Cells.Item(113, 14).Value = -15425.6
# Row 114
$ws.Cells.Item(114, 8).Value = 1411.8572
$ws.Cells.Item(114, 9).Value = 1376.6
$ws.Cells.Item(114, 10).Value = 1500
$ws.Cells.Item(114, 11).Value = 4129.799999999999
$ws.Cells.Item(114, 12).Value = 4500
$ws.Cells.Item(114, 13).Value = -875.7999999999993
$ws.Cells.Item(114, 14).Value = -11008
# Row 122
$ws.Cells.Item(122, 8).Value = 1463.091
$ws.Cells.Item(122, 9).Value = 906.2857
$ws.Cells.Item(122, 11).Value = 8156.571300000001
$ws.Cells.Item(122, 13).Value = -5706.571300000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 1805.1666
$ws.Cells.Item(80, 9).Value = 1626.2
$ws.Cells.Item(80, 10).Value = 2700
$ws.Cells.Item(80, 11).Value = 1626.2
$ws.Cells.Item(80, 12).Value = 2700
$ws.Cells.Item(80, 13).Value = -628.2
$ws.Cells.Item(80, 14).Value = -4696
# Row 83
$ws.Cells.Item(83, 8).Value = 1805.1666
$ws.Cells.Item(83, 9).Value = 1626.2
$ws.Cells.Item(83, 10).Value = 2700
$ws.Cells.Item(83, 11).Value = 8131
$ws.Cells.Item(83, 12).Value = 13500
$ws.Cells.Item(83, 13).Value = -3139
$ws.Cells.Item(83, 14).Value = -23484
# Row 122
$ws.Cells.Item(122, 8).Value = 9999
$ws.Cells.Item(122, 10).Value = 9999
$ws.Cells.Item(122, 12).Value = 29997
$ws.Cells.Item(122, 14).Value = -34897

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Cells.Item(18, 8).Value = 2500
$ws.Cells.Item(18, 9).Value = 2500
$ws.Cells.Item(18, 11).Value = 2500
$ws.Cells.Item(18, 13).Value = -2328
# Row 122
$ws.Cells.Item(122, 8).Value = 4414.65
$ws.Cells.Item(122, 9).Value = 4157.1924
$ws.Cells.Item(122, 10).Value = 4892.7856
$ws.Cells.Item(122, 11).Value = 12471.5772
$ws.Cells.Item(122, 12).Value = 14678.3568
$ws.Cells.Item(122, 13).Value = -10021.5772
$ws.Cells.Item(122, 14).Value = -19578.3568
# Row 131
$ws.Cells.Item(131, 8).Value = 88388.75
$ws.Cells.Item(131, 10).Value = 88388.75
$ws.Cells.Item(131, 12).Value = 88388.75
$ws.Cells.Item(131, 14).Value = -98468.75
# Row 132
$ws.Cells.Item(132, 8).Value = 100003256
$ws.Cells.Item(132, 9).Value = 3076.4614
$ws.Cells.Item(132, 11).Value = 9229.3842
$ws.Cells.Item(132, 13).Value = -6699.3842

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Cells.Item(13, 8).Value = 1624.875
$ws.Cells.Item(13, 9).Value = 1666.5
$ws.Cells.Item(13, 11).Value = 1666.5
$ws.Cells.Item(13, 13).Value = -1526.5
# Row 64
$ws.Cells.Item(64, 8).Value = 40000
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 40000
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 40000
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -40496
# Row 67
$ws.Cells.Item(67, 8).Value = 40000
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 40000
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 40000
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -41716
# Row 107
$ws.Cells.Item(107, 8).Value = 304
$ws.Cells.Item(107, 9).Value = 263.45456
$ws.Cells.Item(107, 11).Value = 790.36368
$ws.Cells.Item(107, 13).Value = 1129.63632
# Row 113
$ws.Cells.Item(113, 8).Value = 794.7778
$ws.Cells.Item(113, 9).Value = 330.625
$ws.Cells.Item(113, 10).Value = 1166.1
$ws.Cells.Item(113, 11).Value = 991.875
$ws.Cells.Item(113, 12).Value = 3498.3
$ws.Cells.Item(113, 13).Value = 1178.125
$ws.Cells.Item(113, 14).Value = -7838.299999999999
# Row 132
$ws.Cells.Item(132, 8).Value = 4341.3823
$ws.Cells.Item(132, 9).Value = 4577.893
$ws.Cells.Item(132, 11).Value = 13733.679
$ws.Cells.Item(132, 13).Value = -11203.679
